$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 15:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 739502
$ws.Range("C4").Value = 710
$ws.Range("D4").Value = 68442
$ws.Range("E4").Value = 632020
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 39040

# Row 9 - Reino Unido
$ws.Range("B9").Value = 120067
$ws.Range("C9").Value = 5850
$ws.Range("E9").Value = 103663
$ws.Range("G9").Value = 596
$ws.Range("H9").Value = 16060

# Row 22 - Austria
$ws.Range("B22").Value = 14696
$ws.Range("C22").Value = 25
$ws.Range("E22").Value = 3752

# Row 42 - Serbia
$ws.Range("F42").Value = 120

# Row 55 - Argentina
$ws.Range("D55").Value = 709
$ws.Range("E55").Value = 1998
$ws.Range("F55").Value = 123

# Row 116 - Sri Lanka
$ws.Range("D116").Value = 96
$ws.Range("E116").Value = 166

# Row 158 - Uganda
$ws.Range("D158").Value = 28
$ws.Range("E158").Value = 27
